# "scrol down and negative test"
# Adds a new "Registration" worksheet after "Address" with a small
# login/registration form (Name / Email / Mobile / Password) and a set of
# rows that each omit exactly one required field (negative tests), plus a
# fully-populated row. The Email column values are hyperlinked (mailto:).

$wb = $excel.ActiveWorkbook
$addressSheet = $wb.Sheets("Address")

# De-select / scroll the Address sheet away from its old selection.
$addressSheet.Range("A1:D1").Select()

# Add the new sheet right after "Address" and name it.
$ws = $wb.Worksheets.Add([System.Type]::Missing, $addressSheet)
$ws.Name = "Registration"

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Mobile"
$ws.Range("D1").Value = "Password"
$ws.Range("A1:D1").Borders.Color = 0
$ws.Range("A1:D1").Interior.Pattern = -4142
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108

# ---- Data rows (each one is a negative test, missing one field) -----
# Row 2: missing Name
$ws.Range("B2").Value = "arun5603@gmail.com"
$ws.Range("C2").Value = "9899001068"
$ws.Range("D2").Value = 12345

# Row 3: missing Email
$ws.Range("A3").Value = "Arun"
$ws.Range("C3").Value = "9899001068"
$ws.Range("D3").Value = 12345

# Row 4: missing Mobile
$ws.Range("A4").Value = "Arun"
$ws.Range("B4").Value = "arun5603@gmail.com"
$ws.Range("D4").Value = 12345

# Row 5: missing Password
$ws.Range("A5").Value = "Arun"
$ws.Range("B5").Value = "arun5603@gmail.com"
$ws.Range("C5").Value = "9899001068"

# Row 6: all fields present (positive / baseline case)
$ws.Range("A6").Value = "Arun"
$ws.Range("B6").Value = "arun5603@gmail.com"
$ws.Range("C6").Value = "9899001068"
$ws.Range("D6").Value = 12345

# ---- Hyperlink the email addresses -----------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:arun5603@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:arun5603@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:arun5603@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:arun5603@gmail.com")

# ---- Borders around the whole table ----------------------------------
$ws.Range("A1:D6").Borders.Color = 0

# ---- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 20.28515625
$ws.Columns.Item(3).ColumnWidth = 12.28515625
$ws.Columns.Item(4).ColumnWidth = 10.7109375

# ---- View / selection state --------------------------------------------
$ws.Range("I10").Select()
$wb.Worksheets("Registration").Activate()
